# Add "90I" / "Runaway" back into the UCROffenseCodeType code table.
#
# A new row is inserted immediately before the existing "910 / 90J / Trespass
# of Real Property" row (row 58), pushing it and everything below it down by
# one row. The new row reuses the same formatting/column layout pattern
# (StateCode/StateDescription/FBICode/FBIDescription = "90I"/"Runaway",
# OffenseCategory1-4 = "Group B"/"Other"/"Group B Offenses (Society)"/"Society"),
# matching the sibling "Group B Offenses (Society)" rows around it (901-910).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCROffenseCodeType")

# Bring the UCROffenseCodeType sheet to the front (it becomes the active tab).
$ws.Activate()

# Insert a new blank row at row 58 - existing row 58 ("910"/"90J"/"Trespass of
# Real Property") and everything below shifts down to row 59+.
$ws.Rows.Item(58).Insert()
$ws.Rows.Item(58).RowHeight = 15

$ws.Range("A58").Value = 909
$ws.Range("B58").Value = "90I"
$ws.Range("C58").Value = "Runaway"
$ws.Range("D58").Value = "90I"
$ws.Range("E58").Value = "Runaway"
$ws.Range("F58").Value = "Group B"
$ws.Range("G58").Value = "Other"
$ws.Range("H58").Value = "Group B Offenses (Society)"
$ws.Range("I58").Value = "Society"

# Reflect the row the editor ended up selecting after inserting the new row.
$ws.Rows.Item(58).Select()
